$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the currency label "ДЕН" to "МКД" across column C (rows 2-23)
for ($r = 2; $r -le 23; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "ДЕН") {
        $cell.Value = "МКД"
    }
}

# D14:D23 lose their distinct (filled) style so they match D2:D13
$ws.Range("D14:D23").Style = $ws.Range("D2").Style

# Update the active selection
$ws.Range("C2").Select()
